# Apply the "Refactor to Fortune-of-the-Day / DynamoDB / CloudFormation" edits
# described by the commit message: update a couple of instructional strings
# and append a new "B. CloudFormation" section (rows 66-72) that mirrors the
# existing "A. Manual" section's structure/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text tweaks to two existing cells
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = "2.5. SSH to EC2, then run jar file"
$ws.Range("E44").Value = "[Reset] button: To reset app state to initialized state (the state on above image) (take some time to delete, create, insert data to DynamoDB)"

# ---------------------------------------------------------------------------
# 2. New "B. CloudFormation" block, rows 66-72 (mirrors rows 1/13/14/18/19/22
#    of the "A. Manual" block above it)
# ---------------------------------------------------------------------------
$ws.Range("A66").Value = "B. CloudFormation"
$ws.Range("C67").Value = "①．Run advanced-fortune-ec2.yml to start EC2 Instance"
$ws.Range("C68").Value = "②． Copy jar file to EC2 instance"
$ws.Range("D69").Value = "Follow the link below (for windows user)"
$ws.Range("E70").Value = "https://intellipaat.com/community/43019/how-to-copy-files-from-local-machine-to-my-aws-instance"
$ws.Range("C71").Value = "③．SSH to EC2, then run jar file"
$ws.Range("E72").Value = "java -jar aws-0.0.1-SNAPSHOT.jar"

# Row height for the new section header (same as row 1's "A. Manual" header)
$ws.Rows.Item(66).RowHeight = $ws.Rows.Item(1).RowHeight

# Copy cell formatting from the analogous cells in the "A. Manual" section
$ws.Range("A1").Copy()
$ws.Range("A66").PasteSpecial(-4122)

$ws.Range("D18").Copy()
$ws.Range("D69").PasteSpecial(-4122)

$ws.Range("E19").Copy()
$ws.Range("E70").PasteSpecial(-4122)

$ws.Range("D22").Copy()
$ws.Range("E72").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Hyperlink for the newly added URL cell (Hyperlinks.Add applies its own
# style, so re-apply the hyperlink-cell formatting afterwards to land on the
# same shared style slot as the other hyperlink cells, e.g. E19)
$ws.Hyperlinks.Add($ws.Range("E70"), "https://intellipaat.com/community/43019/how-to-copy-files-from-local-machine-to-my-aws-instance") | Out-Null

$ws.Range("E19").Copy()
$ws.Range("E70").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the active selection (matches the saved sheetView state)
# ---------------------------------------------------------------------------
$ws.Range("J3").Select()
